$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Fitness (column C) values for rows 2-133 (Generation 0-131 of Run 3),
# replacing the prior log data with the values from the re-run.
$fitnessValues = @(8955, 8955, 8955, 8955, 8955, 8955, 8955, 8955, 8955, 8955, 8955, 8471, 8471, 8471, 8471, 8071, 8071, 8071, 8071, 8071, 8071, 8071, 8003, 8003, 8003, 8003, 7927, 7785, 7785, 7785, 7785, 7736, 7736, 7736, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7320, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295)

for ($i = 0; $i -lt $fitnessValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $fitnessValues[$i]
}
